{"js": "// The document is a single-column table where every row contains one\n// paragraph whose run(s) hold a numeric/text value (some rows pack several\n// tab-separated values into one run). Update the affected rows' text by\n// (0-based) paragraph index, matching the target OOXML from the diff.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"185\",\n  4: \"0.00003\",\n  5: \"0.00013\",\n  6: \"0.00006\",\n  7: \"0.00002\",\n  8: \"0.00006\",\n  9: \"0.00006\",\n  10: \"0.00013\",\n  11: \"0.00761\",\n  43: \"99.99\",\n  44: \"0.01\",\n  45: \"80\",\n};\n\nfor (const idx of Object.keys(updates)) {\n  const i = parseInt(idx, 10);\n  paragraphs.items[i].insertText(updates[idx], \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The document body is a single-column table; each row holds one cell with\n# one paragraph of text (a handful of rows pack several tab-separated\n# values into one run). Update the affected rows' text in place, using\n# 1-based Table.Cell(row, column) addressing.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"185\"\n    5  = \"0.00003\"\n    6  = \"0.00013\"\n    7  = \"0.00006\"\n    8  = \"0.00002\"\n    9  = \"0.00006\"\n    10 = \"0.00006\"\n    11 = \"0.00013\"\n    12 = \"0.00761\"\n    44 = \"99.99\"\n    45 = \"0.01\"\n    46 = \"80\"\n}\n\nforeach ($rowNum in $updates.Keys) {\n    $cell = $t.Cell($rowNum, 1)\n    $cell.Range.Text = $updates[$rowNum]\n}\n"}
